$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-order existing "Installed" rows and interleave the corresponding
# "Not Installed" / "Empty" rows, then append the three new screw-type rows.
$ws.Range("A3").Value = "Empty Panel"
$ws.Range("B3").Value = -1

$ws.Range("A4").Value = "Battery Installed"
$ws.Range("B4").Value = 6

$ws.Range("A5").Value = "Battery Not Installed"
$ws.Range("B5").Value = -1

$ws.Range("A6").Value = "Battery Cushion Installed"
$ws.Range("B6").Value = 9

$ws.Range("A7").Value = "Battery Cover Installed"
$ws.Range("B7").Value = 8

$ws.Range("A8").Value = "Screws Installed"
$ws.Range("B8").Value = 0

$ws.Range("A9").Value = "Screws Not Installed"
$ws.Range("B9").Value = -1

$ws.Range("A10").Value = "U Clamp Installed"
$ws.Range("B10").Value = 1

$ws.Range("A11").Value = "M8 x 35 Screw"
$ws.Range("B11").Value = -1

$ws.Range("A12").Value = "M_F Spacer Screw"
$ws.Range("B12").Value = -1

$ws.Range("A13").Value = "1by4 x 1by2 Screw"
$ws.Range("B13").Value = -1

$ws.Range("A14").Value = "1by4 x 1 Screw"
$ws.Range("B14").Value = -1

# Grow the table so the new rows are included.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:B14"))

# Leftover stray column width (cosmetic, matches original author's session).
$ws.Columns("G:G").ColumnWidth = 22.67

# Match the saved selection/active cell from the edit.
$ws.Range("B11").Select() | Out-Null
